# Jogos_da_Semana_FlashScore_2024-10-25.xlsx update
# - inserts two new "Odd_CS_*_HT" columns (Odd_CS_3-3_HT before the old AW,
#   Odd_CS_4-4_HT after the old BB / new BC) which shifts the existing
#   Odd_CS_0-1_HT..Odd_CS_2-3_HT columns one slot to the right
# - refreshes the single data row (row 2) with the new match's details/odds

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at AW - this shifts the old AW1:BB1 ("Odd_CS_0-1_HT"
#    .. "Odd_CS_2-3_HT") headers (and their row-2 values) one column to the
#    right (to AX1:BC1), preserving their styles/values automatically.
$ws.Range("AW1").EntireColumn.Insert()

# 2) Insert another new column at BD (now just past the old data, which after
#    step 1 ends at BC) so the new cell inherits the same header style as its
#    neighbour without disturbing anything else.
$ws.Range("BD1").EntireColumn.Insert()

# 3) Fill in the two new header cells.
$ws.Range("AW1").Value = "Odd_CS_3-3_HT"
$ws.Range("BD1").Value = "Odd_CS_4-4_HT"

# 4) Update row 2 with the new match's data.
$ws.Range("A2").Value = "QB5xzkQh"
$ws.Range("C2").Value = "05:35"
$ws.Range("D2").Value = "AUSTRALIA - A-LEAGUE"
$ws.Range("E2").Value = "Macarthur FC"
$ws.Range("F2").Value = "Newcastle Jets"
$ws.Range("G2").Value = 1.7
$ws.Range("H2").Value = 4.33
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.2
$ws.Range("K2").Value = 2.6
$ws.Range("L2").Value = 4.33
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 19
$ws.Range("O2").Value = 1.14
$ws.Range("P2").Value = 5.5
$ws.Range("Q2").Value = 1.44
$ws.Range("R2").Value = 2.75
$ws.Range("S2").Value = 1.22
$ws.Range("T2").Value = 4
$ws.Range("U2").Value = 1.5
$ws.Range("V2").Value = 2.5
$ws.Range("W2").Value = 12
$ws.Range("X2").Value = 11
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 12
$ws.Range("AB2").Value = 17
$ws.Range("AC2").Value = 21
$ws.Range("AD2").Value = 9
$ws.Range("AE2").Value = 12
$ws.Range("AF2").Value = 34
$ws.Range("AG2").Value = 101
$ws.Range("AH2").Value = 19
$ws.Range("AI2").Value = 29
$ws.Range("AJ2").Value = 15
$ws.Range("AK2").Value = 51
$ws.Range("AL2").Value = 29
$ws.Range("AM2").Value = 29
$ws.Range("AN2").Value = 4
$ws.Range("AO2").Value = 8.5
$ws.Range("AP2").Value = 15
$ws.Range("AQ2").Value = 23
$ws.Range("AR2").Value = 41
$ws.Range("AS2").Value = 81
$ws.Range("AT2").Value = 4
$ws.Range("AU2").Value = 7
$ws.Range("AV2").Value = 41
$ws.Range("AW2").Value = 301
$ws.Range("AX2").Value = 6.5
$ws.Range("AY2").Value = 21
$ws.Range("AZ2").Value = 21
$ws.Range("BA2").Value = 67
$ws.Range("BB2").Value = 67
$ws.Range("BC2").Value = 101
$ws.Range("BD2").Value = 151

Write-Host "edit.ps1 applied"
